$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.650.51'
$ws.Range("E2").Value = '  -1.97%  '

$ws.Range("D3").Value = '2.004.80'
$ws.Range("E3").Value = '  -4.06%  '

$ws.Range("E4").Value = '  +1.01%  '

$ws.Range("D5").Value = '330.99'

$ws.Range("D6").Value = '1.013'
$ws.Range("E6").Value = '  +0.89%  '

$ws.Range("D7").Value = '0.5019'
$ws.Range("E7").Value = '  -3.64%  '

$ws.Range("D8").Value = '0.4247'
$ws.Range("E8").Value = '  -3.52%  '

$ws.Range("D9").Value = '54.03'
$ws.Range("E9").Value = '  -0.47%  '

$ws.Range("D10").Value = '0.09026'
$ws.Range("E10").Value = '  -3.11%  '

$ws.Range("D11").Value = '1.120'
$ws.Range("E11").Value = '  -4.13%  '

$ws.Range("D12").Value = '23.38'
$ws.Range("E12").Value = '  -5.46%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '8.094'
$ws.Range("E13").Value = '  -6.65%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.980.35'
$ws.Range("E14").Value = '  -6.48%  '

$ws.Range("D15").Value = '6.495'
$ws.Range("E15").Value = '  -5.87%  '

$ws.Range("D16").Value = '1.014'
$ws.Range("E16").Value = '  +0.89%  '

$ws.Range("D17").Value = '94.23'
$ws.Range("E17").Value = '  -6.87%  '

$ws.Range("D18").Value = '0.00001118'
$ws.Range("E18").Value = '  -3.28%  '

$ws.Range("D19").Value = '0.06676'
$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").Value = '19.80'
$ws.Range("E20").Value = '  -6.49%  '

$ws.Range("D21").Value = '1.013'
$ws.Range("E21").Value = '  +0.81%  '

$ws.Range("D22").Value = '5.984'
$ws.Range("E22").Value = '  -6.04%  '

$ws.Range("D23").Value = '29.654.16'
$ws.Range("E23").Value = '  -2.02%  '

$ws.Range("D24").Value = '12.03'
$ws.Range("E24").Value = '  -3.84%  '

$ws.Range("D25").Value = '2.289'
$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("D26").Value = '159.67'
$ws.Range("E26").Value = '  -1.46%  '

$ws.Range("D27").Value = '20.76'
$ws.Range("E27").Value = '  -4.29%  '

$ws.Range("D28").Value = '6.438'
$ws.Range("E28").Value = '  -3.70%  '

$ws.Range("D29").Value = '2.304'
$ws.Range("E29").Value = '  -8.42%  '

$ws.Range("D30").Value = '128.98'
$ws.Range("E30").Value = '  -2.96%  '

$ws.Range("D31").Value = '1.056'
$ws.Range("E31").Value = '  -6.41%  '

$ws.Range("D32").Value = '0.09946'
$ws.Range("E32").Value = '  -4.90%  '

$ws.Range("D33").Value = '1.575'
$ws.Range("E33").Value = '  -5.56%  '

$ws.Range("D34").Value = '5.847'
$ws.Range("E34").Value = '  -5.86%  '

$ws.Range("D35").Value = '3.806'
$ws.Range("E35").Value = '  -1.41%  '

$ws.Range("D36").Value = '0.02475'
$ws.Range("E36").Value = '  -5.80%  '

$ws.Range("D37").Value = '9.417'
$ws.Range("E37").Value = '  -7.65%  '

$ws.Range("D38").Value = '1.311'
$ws.Range("E38").Value = '  -2.91%  '

$ws.Range("D39").Value = '0.06365'
$ws.Range("E39").Value = '  -5.67%  '

$ws.Range("D40").Value = '0.6591'
$ws.Range("E40").Value = '  -5.54%  '

$ws.Range("D41").Value = '11.71'
$ws.Range("E41").Value = '  -6.48%  '

$ws.Range("D42").Value = '0.2055'
$ws.Range("E42").Value = '  -7.09%  '

$ws.Range("D43").Value = '1.012'
$ws.Range("E43").Value = '  +0.76%  '

$ws.Range("D44").Value = '0.6354'
$ws.Range("E44").Value = '  -6.99%  '

$ws.Range("D45").Value = '13.46'
$ws.Range("E45").Value = '  -5.80%  '

$ws.Range("D46").Value = '2.211'
$ws.Range("E46").Value = '  -5.51%  '

$ws.Range("D47").Value = '1.297'
$ws.Range("E47").Value = '  -4.65%  '

$ws.Range("D48").Value = '3.528'
$ws.Range("E48").Value = '  -2.93%  '

$ws.Range("D49").Value = '0.00000000339'
$ws.Range("E49").Value = '  -2.05%  '

$ws.Range("D50").Value = '0.07004'
$ws.Range("E50").Value = '  -3.00%  '

$ws.Range("D51").Value = '1.126'
$ws.Range("E51").Value = '  -7.47%  '
